$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.513.88'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '1.869.62'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").Value = "'312.29"
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").Value = "'0.4782"
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = "'0.3779"
$ws.Range("E8").Value = '  +2.98%  '
$ws.Range("D9").Value = "'0.07349"
$ws.Range("E9").Value = '  +1.79%  '
$ws.Range("D10").Value = "'0.9370"
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").Value = "'20.70"
$ws.Range("E11").Value = '  +4.90%  '
$ws.Range("D12").Value = "'0.07848"
$ws.Range("E12").Value = '  +2.00%  '
$ws.Range("D13").Value = '1.856.45'
$ws.Range("E13").Value = '  -1.98%  '
$ws.Range("D14").Value = "'5.443"
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").Value = "'6.570"
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("D16").Value = "'90.74"
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("D17").Value = "'1.015"
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Value = "'0.000008897"
$ws.Range("E18").Value = '  +2.97%  '
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").Value = "'14.92"
$ws.Range("E20").Value = '  +2.48%  '
$ws.Range("D21").Value = '27.508.68'
$ws.Range("E21").Value = '  +1.95%  '
$ws.Range("D22").Value = "'5.128"
$ws.Range("E22").Value = '  +1.44%  '
$ws.Range("D23").Value = "'10.72"
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").Value = "'1.955"
$ws.Range("E24").Value = '  +1.62%  '
$ws.Range("D25").Value = "'153.83"
$ws.Range("E26").Value = '  +1.88%  '
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("D28").Value = "'115.81"
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("D29").Value = "'4.993"
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("D30").Value = "'0.08923"
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("D32").Value = "'1.215"
$ws.Range("E32").Value = '  +3.35%  '
$ws.Range("D33").Value = "'4.611"
$ws.Range("E33").Value = '  +2.65%  '
$ws.Range("D34").Value = "'0.7542"
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("D35").Value = "'2.710"
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("D36").Value = "'0.02049"
$ws.Range("E36").Value = '  +4.57%  '
$ws.Range("D37").Value = "'1.118"
$ws.Range("E37").Value = '  +0.92%  '
$ws.Range("D38").Value = "'3.003"
$ws.Range("E38").Value = '  +0.94%  '
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("E40").Value = '  +2.90%  '
$ws.Range("E41").Value = '  +1.73%  '
$ws.Range("E42").Value = '  +1.08%  '
$ws.Range("D43").Value = "'8.496"
$ws.Range("E43").Value = '  +3.59%  '
$ws.Range("D44").Value = "'10.67"
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("D45").Value = "'0.4811"
$ws.Range("E45").Value = '  +1.81%  '
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("E47").Value = '  +3.61%  '
$ws.Range("D48").Value = "'102.83"
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("D49").Value = "'67.43"
$ws.Range("E49").Value = '  +1.96%  '
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").Value = "'0.9237"
$ws.Range("E51").Value = '  +4.31%  '
